$wb = $excel.ActiveWorkbook

# --- Step 1: add a "2022-Q1" sheet, right before "总计" -------------------
# Duplicate "2021-Q4" (same column layout/headers/styling as the sheet we
# need) and slot the copy in immediately before "总计".
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$q4.Copy($total)

# Re-resolve positionally: Copy() inserts the duplicate immediately before
# "总计", and any object/index references captured before the Copy() call
# are stale afterwards, so look the new sheet up fresh by its position.
$total = $wb.Worksheets.Item("总计")
$newQ1 = $wb.Worksheets.Item($total.Index - 1)
$newQ1.Name = "2022-Q1"

# Overwrite the fund rows with the 2022-Q1 figures (text-like numeric
# columns keep their original text formatting via a leading apostrophe).
$newQ1.Range("D2").Value = "'0.20"
$newQ1.Range("E2").Value = "'93.65"
$newQ1.Range("F2").Value = "'4.21"
$newQ1.Range("G2").Value = "'0.0084"
$newQ1.Range("H2").Value = 4

$newQ1.Range("D3").Value = "'0.06"
$newQ1.Range("E3").Value = "'93.65"
$newQ1.Range("F3").Value = "'4.21"
$newQ1.Range("G3").Value = "'0.0025"
$newQ1.Range("H3").Value = 4

# --- Step 2: add the 2022-Q1 row to "总计" ---------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

# The blank row Insert() picked up borders from the row beneath it; strip
# that back off the data cells and restore the index-column styling (taken
# from the sibling index cell) so A2 matches A3:A5.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Renumber the index column for the rows that shifted down one place.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

$wb.Worksheets.Item("2021-Q2").Select()
